$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2416.6667
$ws.Range("J2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("N2").Value = -5226
$ws.Range("H28").Value = 651.0952
$ws.Range("I28").Value = 671.2105
$ws.Range("K28").Value = 671.2105
$ws.Range("M28").Value = -186.2105
$ws.Range("H51").Value = 4712.857
$ws.Range("I51").Value = 3999
$ws.Range("J51").Value = 4831.8335
$ws.Range("K51").Value = 3999
$ws.Range("L51").Value = 4831.8335
$ws.Range("M51").Value = -3515
$ws.Range("N51").Value = -5799.8335
$ws.Range("H70").Value = 4383.25
$ws.Range("I70").Value = 6110.8
$ws.Range("J70").Value = 3149.2856
$ws.Range("K70").Value = 18332.4
$ws.Range("L70").Value = 9447.856800000001
$ws.Range("M70").Value = -18062.4
$ws.Range("N70").Value = -9987.856800000001
$ws.Range("H73").Value = 4383.25
$ws.Range("I73").Value = 6110.8
$ws.Range("J73").Value = 3149.2856
$ws.Range("K73").Value = 18332.4
$ws.Range("L73").Value = 9447.856800000001
$ws.Range("M73").Value = -17396.4
$ws.Range("N73").Value = -11319.8568

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22631.773
$ws.Range("I32").Value = 19899.953
$ws.Range("K32").Value = 19899.953
$ws.Range("M32").Value = -19612.953
$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -35976
$ws.Range("H55").Value = 24996.428
$ws.Range("J55").Value = 24996.428
$ws.Range("L55").Value = 24996.428
$ws.Range("N55").Value = -25626.428
$ws.Range("H132").Value = 8070.143
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 26350.75
$ws.Range("J82").Value = 39998.57
$ws.Range("L82").Value = 39998.57
$ws.Range("N82").Value = -40764.57
$ws.Range("H85").Value = 26350.75
$ws.Range("J85").Value = 39998.57
$ws.Range("L85").Value = 39998.57
$ws.Range("N85").Value = -42650.57

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5500
$ws.Range("I16").Value = 5500
$ws.Range("K16").Value = 5500
$ws.Range("M16").Value = -5213
$ws.Range("H31").Value = 5742.8
$ws.Range("I31").Value = 4949.5
$ws.Range("J31").Value = 6271.6665
$ws.Range("K31").Value = 4949.5
$ws.Range("L31").Value = 6271.6665
$ws.Range("M31").Value = -4654.5
$ws.Range("N31").Value = -6861.6665
$ws.Range("H34").Value = 5742.8
$ws.Range("I34").Value = 4949.5
$ws.Range("J34").Value = 6271.6665
$ws.Range("K34").Value = 4949.5
$ws.Range("L34").Value = 6271.6665
$ws.Range("M34").Value = -4747.5
$ws.Range("N34").Value = -6675.6665
$ws.Range("H41").Value = 18885.715
$ws.Range("I41").Value = 11100
$ws.Range("J41").Value = 22000
$ws.Range("K41").Value = 11100
$ws.Range("L41").Value = 22000
$ws.Range("M41").Value = -10672
$ws.Range("N41").Value = -22856
$ws.Range("H74").Value = 33438
$ws.Range("J74").Value = 33438
$ws.Range("L74").Value = 33438
$ws.Range("N74").Value = -35186
$ws.Range("H77").Value = 33438
$ws.Range("J77").Value = 33438
$ws.Range("L77").Value = 100314
$ws.Range("N77").Value = -109050
$ws.Range("H94").Value = 2032.1666
$ws.Range("J94").Value = 1648.5
$ws.Range("L94").Value = 1648.5
$ws.Range("N94").Value = -2550.5
$ws.Range("H113").Value = 5500
$ws.Range("I113").Value = 5500
$ws.Range("K113").Value = 5500
$ws.Range("M113").Value = -3330
$ws.Range("H132").Value = 1532.6666
$ws.Range("I132").Value = 1399.5714
$ws.Range("K132").Value = 4198.7142
$ws.Range("M132").Value = -1668.7142
$ws.Range("H134").Value = 8676.5
$ws.Range("I134").Value = 7096
$ws.Range("J134").Value = 14998.5
$ws.Range("K134").Value = 21288
$ws.Range("L134").Value = 44995.5
$ws.Range("M134").Value = -18753
$ws.Range("N134").Value = -50065.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 400.5
$ws.Range("I2").Value = 332.33334
$ws.Range("K2").Value = 1994.00004
$ws.Range("M2").Value = -1881.00004
$ws.Range("H38").Value = 119.72727
$ws.Range("I38").Value = 127
$ws.Range("K38").Value = 381
$ws.Range("M38").Value = -34
$ws.Range("H40").Value = 79.40000000000001
$ws.Range("J40").Value = 39.4
$ws.Range("L40").Value = 157.6
$ws.Range("N40").Value = -295.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 33400
$ws.Range("I62").Value = 27000
$ws.Range("J62").Value = 35000
$ws.Range("K62").Value = 27000
$ws.Range("L62").Value = 35000
$ws.Range("M62").Value = -26314
$ws.Range("N62").Value = -36372
$ws.Range("H65").Value = 33400
$ws.Range("I65").Value = 27000
$ws.Range("J65").Value = 35000
$ws.Range("K65").Value = 81000
$ws.Range("L65").Value = 105000
$ws.Range("M65").Value = -77568
$ws.Range("N65").Value = -111864
$ws.Range("H102").Value = 1901.4333
$ws.Range("I102").Value = 1811.7916
$ws.Range("K102").Value = 1811.7916
$ws.Range("M102").Value = -189.7916
$ws.Range("H109").Value = 90000
$ws.Range("J109").Value = 90000
$ws.Range("L109").Value = 90000
$ws.Range("N109").Value = -92080
$ws.Range("H122").Value = 35307.6
$ws.Range("I122").Value = 36915.562
$ws.Range("K122").Value = 110746.686
$ws.Range("M122").Value = -108296.686
$ws.Range("H132").Value = 1450
$ws.Range("I132").Value = 900
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 2700
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -170
$ws.Range("N132").Value = -11060
$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2186.125
$ws.Range("I46").Value = 2186.125
$ws.Range("K46").Value = 2186.125
$ws.Range("M46").Value = -1998.125
$ws.Range("H93").Value = 2246.125
$ws.Range("I93").Value = 1700
$ws.Range("K93").Value = 1700
$ws.Range("M93").Value = -452

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 23999.666
$ws.Range("I39").Value = 24499.5
$ws.Range("J39").Value = 23000
$ws.Range("K39").Value = 24499.5
$ws.Range("L39").Value = 23000
$ws.Range("M39").Value = -24086.5
$ws.Range("N39").Value = -23826
$ws.Range("H101").Value = 16338
$ws.Range("J101").Value = 16338
$ws.Range("L101").Value = 16338
$ws.Range("N101").Value = -22828
$ws.Range("H126").Value = 4088.1
$ws.Range("I126").Value = 3610.25
$ws.Range("K126").Value = 10830.75
$ws.Range("M126").Value = -8360.75
